$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = [double]"22.34000000000005"
$ws.Range("H2").Value = [double]"1.099376145674569e-11"
$ws.Range("I2").Value = [double]"1.099376145674569e-11"
$ws.Range("L2").Value = [double]"63.0751608361858"
$ws.Range("M2").Value = "[48.849320645998105, 77.3010010263735]"
$ws.Range("N2").Value = [double]"1.590239051552089e-11"
$ws.Range("O2").Value = [double]"1.590239051552089e-11"
$ws.Range("P2").Value = [double]"1.641552918091964"
$ws.Range("Q2").Value = "[1.415131825941348, 1.8679740102425804]"
$ws.Range("R2").Value = [double]"0"
$ws.Range("S2").Value = [double]"0"
$ws.Range("T2").Value = [double]"62.69549740817033"
$ws.Range("U2").Value = "[53.64355882976873, 71.74743598657193]"
$ws.Range("V2").Value = [double]"0"
$ws.Range("W2").Value = [double]"0"
$ws.Range("X2").Value = [double]"16.50342342342346"
$ws.Range("Y2").Value = [double]"15.69837837837841"
$ws.Range("Z2").Value = [double]"17.30846846846851"

# Row 3
$ws.Range("F3").Value = [double]"22.34000000000005"
$ws.Range("H3").Value = [double]"6.155486897974072e-09"
$ws.Range("I3").Value = [double]"6.155486897974072e-09"
$ws.Range("L3").Value = [double]"57.05317747093013"
$ws.Range("M3").Value = "[38.6618786805533, 75.44447626130696]"
$ws.Range("N3").Value = [double]"1.335941548497033e-07"
$ws.Range("O3").Value = [double]"1.335941548497033e-07"
$ws.Range("P3").Value = [double]"1.893131909370426"
$ws.Range("Q3").Value = "[1.528342372016657, 2.257921446724196]"
$ws.Range("R3").Value = [double]"1.27675647831893e-13"
$ws.Range("S3").Value = [double]"1.27675647831893e-13"
$ws.Range("T3").Value = [double]"57.74415297076485"
$ws.Range("U3").Value = "[47.36416925780092, 68.12413668372878]"
$ws.Range("V3").Value = [double]"1.310063169057685e-14"
$ws.Range("W3").Value = [double]"1.310063169057685e-14"
$ws.Range("X3").Value = [double]"15.60892892892896"
$ws.Range("Y3").Value = [double]"14.31191191191195"
$ws.Range("Z3").Value = [double]"16.90594594594598"

# Row 4
$ws.Range("F4").Value = [double]"22.34000000000005"
$ws.Range("H4").Value = [double]"1.600923281719346e-07"
$ws.Range("I4").Value = [double]"1.600923281719346e-07"
$ws.Range("L4").Value = [double]"47.74748100188325"
$ws.Range("M4").Value = "[27.453008192853773, 68.04195381091273]"
$ws.Range("N4").Value = [double]"2.18443779922417e-05"
$ws.Range("O4").Value = [double]"2.18443779922417e-05"
$ws.Range("P4").Value = [double]"2.459184639746966"
$ws.Range("Q4").Value = "[2.0440793041375036, 2.8742899753564277]"
$ws.Range("R4").Value = [double]"1.554312234475219e-15"
$ws.Range("S4").Value = [double]"1.554312234475219e-15"
$ws.Range("T4").Value = [double]"48.96312447498212"
$ws.Range("U4").Value = "[38.671128325648375, 59.25512062431587]"
$ws.Range("V4").Value = [double]"1.947775274402375e-12"
$ws.Range("W4").Value = [double]"1.947775274402375e-12"
$ws.Range("X4").Value = [double]"13.59631631631635"
$ws.Range("Y4").Value = [double]"12.12040040040043"
$ws.Range("Z4").Value = [double]"15.07223223223227"

# Row 5
$ws.Range("F5").Value = [double]"22.34000000000005"
$ws.Range("H5").Value = [double]"5.208644443521759e-05"
$ws.Range("I5").Value = [double]"5.208644443521759e-05"
$ws.Range("L5").Value = [double]"42.1479599814714"
$ws.Range("M5").Value = "[21.366522345885663, 62.92939761705714]"
$ws.Range("N5").Value = [double]"0.000178920427601259"
$ws.Range("O5").Value = [double]"0.000178920427601259"
$ws.Range("P5").Value = [double]"3.012658420559582"
$ws.Range("Q5").Value = "[2.408868841491273, 3.6164479996278907]"
$ws.Range("R5").Value = [double]"4.454214774796128e-13"
$ws.Range("S5").Value = [double]"4.454214774796128e-13"
$ws.Range("T5").Value = [double]"58.22235438480802"
$ws.Range("U5").Value = "[45.762259163077175, 70.68244960653887]"
$ws.Range("V5").Value = [double]"3.358868738700949e-12"
$ws.Range("W5").Value = [double]"3.358868738700949e-12"
$ws.Range("X5").Value = [double]"11.62842842842845"
$ws.Range("Y5").Value = [double]"9.48164164164166"
$ws.Range("Z5").Value = [double]"13.77521521521525"

# Row 6
$ws.Range("F6").Value = [double]"24.77000000000043"
$ws.Range("H6").Value = [double]"8.640332671561168e-08"
$ws.Range("I6").Value = [double]"8.640332671561168e-08"
$ws.Range("L6").Value = [double]"46.7228571343253"
$ws.Range("M6").Value = "[31.107592072489844, 62.33812219616076]"
$ws.Range("N6").Value = [double]"2.852413820697564e-07"
$ws.Range("O6").Value = [double]"2.852413820697564e-07"
$ws.Range("P6").Value = [double]"-3.00636894577762"
$ws.Range("Q6").Value = "[-3.3963163822592355, -2.6164215092960044]"
$ws.Range("R6").Value = [double]"0"
$ws.Range("S6").Value = [double]"0"
$ws.Range("T6").Value = [double]"47.52353662405681"
$ws.Range("U6").Value = "[38.05126810044956, 56.99580514766406]"
$ws.Range("V6").Value = [double]"3.745892485085278e-13"
$ws.Range("W6").Value = [double]"3.745892485085278e-13"
$ws.Range("X6").Value = [double]"11.85191191191212"
$ws.Range("Y6").Value = [double]"10.31463463463482"
$ws.Range("Z6").Value = [double]"13.38918918918942"

# Row 7
$ws.Range("F7").Value = [double]"24.77000000000043"
$ws.Range("H7").Value = [double]"1.256064663390788e-08"
$ws.Range("I7").Value = [double]"1.256064663390788e-08"
$ws.Range("L7").Value = [double]"56.24369345565079"
$ws.Range("M7").Value = "[38.936444707376396, 73.55094220392519]"
$ws.Range("N7").Value = [double]"4.827779220839545e-08"
$ws.Range("O7").Value = [double]"4.827779220839545e-08"
$ws.Range("P7").Value = [double]"3.037816319687427"
$ws.Range("Q7").Value = "[2.698184681461504, 3.377447957913351]"
$ws.Range("T7").Value = [double]"60.86224541748998"
$ws.Range("U7").Value = "[50.309285288924706, 71.41520554605526]"
$ws.Range("V7").Value = [double]"3.774758283725532e-15"
$ws.Range("W7").Value = [double]"3.774758283725532e-15"
$ws.Range("X7").Value = [double]"12.79411411411434"
$ws.Range("Y7").Value = [double]"11.4551951951954"
$ws.Range("Z7").Value = [double]"14.13303303303328"

# Row 8
$ws.Range("F8").Value = [double]"24.77000000000043"
$ws.Range("H8").Value = [double]"1.002115082360255e-06"
$ws.Range("I8").Value = [double]"1.002115082360255e-06"
$ws.Range("L8").Value = [double]"51.17428809211368"
$ws.Range("M8").Value = "[31.347664574420364, 71.000911609807]"
$ws.Range("N8").Value = [double]"4.735798172283978e-06"
$ws.Range("O8").Value = [double]"4.735798172283978e-06"
$ws.Range("P8").Value = [double]"2.937184723176043"
$ws.Range("Q8").Value = "[2.4717635893108882, 3.4026058570411974]"
$ws.Range("R8").Value = [double]"2.220446049250313e-16"
$ws.Range("S8").Value = [double]"2.220446049250313e-16"
$ws.Range("T8").Value = [double]"54.21878257551738"
$ws.Range("U8").Value = "[42.56589103981702, 65.87167411121774]"
$ws.Range("V8").Value = [double]"3.818945160105613e-12"
$ws.Range("W8").Value = [double]"3.818945160105613e-12"
$ws.Range("X8").Value = [double]"13.19083083083106"
$ws.Range("Y8").Value = [double]"11.35601601601621"
$ws.Range("Z8").Value = [double]"15.02564564564591"

# Row 9
$ws.Range("F9").Value = [double]"24.77000000000043"
$ws.Range("H9").Value = [double]"3.125646630408596e-08"
$ws.Range("I9").Value = [double]"3.125646630408596e-08"
$ws.Range("L9").Value = [double]"49.81561436897812"
$ws.Range("M9").Value = "[31.786015353585597, 67.84521338437065]"
$ws.Range("N9").Value = [double]"1.374763645367594e-06"
$ws.Range("O9").Value = [double]"1.374763645367594e-06"
$ws.Range("P9").Value = [double]"2.72334258058935"
$ws.Range("Q9").Value = "[2.333395144107734, 3.113290017070966]"
$ws.Range("T9").Value = [double]"48.70309027359311"
$ws.Range("U9").Value = "[38.92990288219251, 58.4762776649937]"
$ws.Range("V9").Value = [double]"4.631850458736153e-13"
$ws.Range("W9").Value = [double]"4.631850458736153e-13"
$ws.Range("X9").Value = [double]"14.0338538538541"
$ws.Range("Y9").Value = [double]"12.49657657657679"
$ws.Range("Z9").Value = [double]"15.5711311311314"

# Row 10
$ws.Range("F10").Value = [double]"24.77000000000043"
$ws.Range("H10").Value = [double]"1.194785714808688e-09"
$ws.Range("I10").Value = [double]"1.194785714808688e-09"
$ws.Range("L10").Value = [double]"56.21156792200861"
$ws.Range("M10").Value = "[36.94970700226878, 75.47342884174843]"
$ws.Range("N10").Value = [double]"4.740944825609716e-07"
$ws.Range("O10").Value = [double]"4.740944825609716e-07"
$ws.Range("P10").Value = [double]"2.232763547596349"
$ws.Range("Q10").Value = "[1.9057108589343477, 2.5598162362583503]"
$ws.Range("T10").Value = [double]"51.94267022017806"
$ws.Range("U10").Value = "[42.09099516073857, 61.794345279617545]"
$ws.Range("V10").Value = [double]"7.66053886991358e-14"
$ws.Range("W10").Value = [double]"7.66053886991358e-14"
$ws.Range("X10").Value = [double]"15.96784784784813"
$ws.Range("Y10").Value = [double]"14.67851851851877"
$ws.Range("Z10").Value = [double]"17.25717717717748"

# Row 11
$ws.Range("F11").Value = [double]"24.77000000000043"
$ws.Range("H11").Value = [double]"4.631357540296754e-05"
$ws.Range("I11").Value = [double]"4.631357540296754e-05"
$ws.Range("L11").Value = [double]"48.74192743814358"
$ws.Range("M11").Value = "[21.41644541767171, 76.06740945861546]"
$ws.Range("N11").Value = [double]"0.0008068716628921724"
$ws.Range("O11").Value = [double]"0.0008068716628921724"
$ws.Range("P11").Value = [double]"2.446605690183042"
$ws.Range("Q11").Value = "[1.9308687580621955, 2.962342622303889]"
$ws.Range("R11").Value = [double]"2.124078690712849e-12"
$ws.Range("S11").Value = [double]"2.124078690712849e-12"
$ws.Range("T11").Value = [double]"63.28468118345029"
$ws.Range("U11").Value = "[49.396872160257544, 77.17249020664303]"
$ws.Range("V11").Value = [double]"7.116751632452178e-12"
$ws.Range("W11").Value = [double]"7.116751632452178e-12"
$ws.Range("X11").Value = [double]"15.12482482482509"
$ws.Range("Y11").Value = [double]"13.09165165165188"
$ws.Range("Z11").Value = [double]"17.1579979979983"

# Row 12
$ws.Range("F12").Value = [double]"24.77000000000043"
$ws.Range("H12").Value = [double]"3.248421087675979e-09"
$ws.Range("I12").Value = [double]"3.248421087675979e-09"
$ws.Range("L12").Value = [double]"51.67145350644154"
$ws.Range("M12").Value = "[34.128434534493366, 69.21447247838972]"
$ws.Range("N12").Value = [double]"3.934068810718117e-07"
$ws.Range("O12").Value = [double]"3.934068810718117e-07"
$ws.Range("P12").Value = [double]"1.956026657190042"
$ws.Range("Q12").Value = "[1.603816069400195, 2.308237244979888]"
$ws.Range("R12").Value = [double]"1.398881011027697e-14"
$ws.Range("S12").Value = [double]"1.398881011027697e-14"
$ws.Range("T12").Value = [double]"58.22123472598854"
$ws.Range("U12").Value = "[48.75262522161407, 67.68984423036302]"
$ws.Range("V12").Value = [double]"4.440892098500626e-16"
$ws.Range("W12").Value = [double]"4.440892098500626e-16"
$ws.Range("X12").Value = [double]"17.05881881881912"
$ws.Range("Y12").Value = [double]"15.67031031031058"
$ws.Range("Z12").Value = [double]"18.44732732732765"

# Row 13
$ws.Range("F13").Value = [double]"24.77000000000043"
$ws.Range("H13").Value = [double]"7.615129549165545e-08"
$ws.Range("I13").Value = [double]"7.615129549165545e-08"
$ws.Range("L13").Value = [double]"56.97758232037678"
$ws.Range("M13").Value = "[36.87813611489189, 77.07702852586168]"
$ws.Range("N13").Value = [double]"8.410270864178671e-07"
$ws.Range("O13").Value = [double]"8.410270864178671e-07"
$ws.Range("P13").Value = [double]"1.767342413731195"
$ws.Range("Q13").Value = "[1.3522370781217337, 2.182447749340657]"
$ws.Range("R13").Value = [double]"5.092415378271653e-11"
$ws.Range("S13").Value = [double]"5.092415378271653e-11"
$ws.Range("T13").Value = [double]"57.11927290355061"
$ws.Range("U13").Value = "[45.24851105532488, 68.99003475177634]"
$ws.Range("V13").Value = [double]"1.375566327510569e-12"
$ws.Range("W13").Value = [double]"1.375566327510569e-12"
$ws.Range("X13").Value = [double]"17.80266266266297"
$ws.Range("Y13").Value = [double]"16.16620620620649"
$ws.Range("Z13").Value = [double]"19.43911911911945"

# Row 14
$ws.Range("F14").Value = [double]"24.77000000000043"
$ws.Range("H14").Value = [double]"6.916254446931447e-08"
$ws.Range("I14").Value = [double]"6.916254446931447e-08"
$ws.Range("L14").Value = [double]"56.75172412104313"
$ws.Range("M14").Value = "[36.187668591697204, 77.31577965038905]"
$ws.Range("N14").Value = [double]"1.405514270080843e-06"
$ws.Range("O14").Value = [double]"1.405514270080843e-06"
$ws.Range("P14").Value = [double]"1.842816111114733"
$ws.Range("Q14").Value = "[1.427710775505271, 2.257921446724196]"
$ws.Range("R14").Value = [double]"1.533373428230789e-11"
$ws.Range("S14").Value = [double]"1.533373428230789e-11"
$ws.Range("T14").Value = [double]"59.79739330031348"
$ws.Range("U14").Value = "[48.02181242092027, 71.57297417970669]"
$ws.Range("V14").Value = [double]"2.55573340268711e-13"
$ws.Range("W14").Value = [double]"2.55573340268711e-13"
$ws.Range("X14").Value = [double]"17.50512512512543"
$ws.Range("Y14").Value = [double]"15.86866866866895"
$ws.Range("Z14").Value = [double]"19.14158158158192"
